$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: 'Datos actualizados a 19 de Mayo de 2020 a las 14:05' -> 'Datos actualizados a 19 de Mayo de 2020 a las 14:35'
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 19 de Mayo de 2020 a las 14:35"

# Row 27: 'Suiza' -> 'Suecia'
$ws.Cells.Item(27, 1).Value = "Suecia"
$ws.Cells.Item(27, 2).Value = 30799
$ws.Cells.Item(27, 3).Value = 422
$ws.Cells.Item(27, 4).Value = 4971
$ws.Cells.Item(27, 5).Value = 22085
$ws.Cells.Item(27, 7).Value = 45
$ws.Cells.Item(27, 8).Value = 3743

# Row 28: 'Suecia' -> 'Suiza'
$ws.Cells.Item(28, 1).Value = "Suiza"
$ws.Cells.Item(28, 2).Value = 30618
$ws.Cells.Item(28, 3).Value = 21
$ws.Cells.Item(28, 4).Value = 27600
$ws.Cells.Item(28, 5).Value = 1132
$ws.Cells.Item(28, 8).Value = 1886

# Row 38: 'Israel' -> 'Kuwait'
$ws.Cells.Item(38, 1).Value = "Kuwait"
$ws.Cells.Item(38, 2).Value = 16764
$ws.Cells.Item(38, 3).Value = 1073
$ws.Cells.Item(38, 4).Value = 4681
$ws.Cells.Item(38, 5).Value = 11962
$ws.Cells.Item(38, 7).Value = 3
$ws.Cells.Item(38, 8).Value = 121

# Row 39: 'Sudafrica' -> 'Israel'
$ws.Cells.Item(39, 1).Value = "Israel"
$ws.Cells.Item(39, 2).Value = 16650
$ws.Cells.Item(39, 3).Value = 7
$ws.Cells.Item(39, 4).Value = 13299
$ws.Cells.Item(39, 5).Value = 3074
$ws.Cells.Item(39, 7).Value = 1
$ws.Cells.Item(39, 8).Value = 277

# Row 40: 'Austria' -> 'Sudafrica'
$ws.Cells.Item(40, 1).Value = "Sudafrica"
$ws.Cells.Item(40, 2).Value = 16433
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = 7298
$ws.Cells.Item(40, 5).Value = 8849
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 286

# Row 41: 'Japon' -> 'Austria'
$ws.Cells.Item(41, 1).Value = "Austria"
$ws.Cells.Item(41, 2).Value = 16321
$ws.Cells.Item(41, 3).Value = 52
$ws.Cells.Item(41, 4).Value = 14678
$ws.Cells.Item(41, 5).Value = 1011
$ws.Cells.Item(41, 7).Value = 3
$ws.Cells.Item(41, 8).Value = 632

# Row 42: 'Colombia' -> 'Japon'
$ws.Cells.Item(42, 1).Value = "Japon"
$ws.Cells.Item(42, 2).Value = 16305
$ws.Cells.Item(42, 4).Value = 11564
$ws.Cells.Item(42, 5).Value = 3992
$ws.Cells.Item(42, 8).Value = 749

# Row 43: 'Kuwait' -> 'Colombia'
$ws.Cells.Item(43, 1).Value = "Colombia"
$ws.Cells.Item(43, 2).Value = 16295
$ws.Cells.Item(43, 4).Value = 3903
$ws.Cells.Item(43, 5).Value = 11800
$ws.Cells.Item(43, 8).Value = 592

# Row 48: 'Dinamarca' -> 'Dinamarca'
$ws.Cells.Item(48, 4).Value = 9416
$ws.Cells.Item(48, 5).Value = 1077
$ws.Cells.Item(48, 7).Value = 3
$ws.Cells.Item(48, 8).Value = 551

# Row 82: 'Croacia' -> 'Croacia'
$ws.Cells.Item(82, 2).Value = 2232
$ws.Cells.Item(82, 3).Value = 4
$ws.Cells.Item(82, 4).Value = 1967
$ws.Cells.Item(82, 5).Value = 169
$ws.Cells.Item(82, 7).Value = 1
$ws.Cells.Item(82, 8).Value = 96

# Row 196: 'Nueva Caledonia' -> 'Belice'
$ws.Cells.Item(196, 1).Value = "Belice"
$ws.Cells.Item(196, 4).Value = 16
$ws.Cells.Item(196, 8).Value = 2

# Row 197: 'Belice' -> 'Nueva Caledonia'
$ws.Cells.Item(197, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(197, 4).Value = 18
$ws.Cells.Item(197, 8).Value = 0

# Row 209: 'Seychelles' -> 'Groenlandia'
$ws.Cells.Item(209, 1).Value = "Groenlandia"

# Row 210: 'Groenlandia' -> 'Montserrat'
$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1

# Row 211: 'Montserrat' -> 'Seychelles'
$ws.Cells.Item(211, 1).Value = "Seychelles"
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0

# Row 215: 'San Bartolome' -> 'Bonaire, San Eustaquio y Saba'
$ws.Cells.Item(215, 1).Value = "Bonaire, San Eustaquio y Saba"

# Row 216: 'Bonaire, San Eustaquio y Saba' -> 'San Bartolome'
$ws.Cells.Item(216, 1).Value = "San Bartolome"
